# Adds the MeanSquaredError family of indices to the list of distance/
# similarity measures on Hoja1 (rows 44-50, right after JensenDifference).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newMeasures = @(
    "MeanSquaredError",
    "RootMeanSquaredError",
    "MeanAbsoluteError",
    "RelativeSquaredError",
    "RootRelativeSquaredError",
    "RelativeAbsoluteError",
    "CorrelationCoefficient"
)

$startRow = 44
for ($i = 0; $i -lt $newMeasures.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newMeasures[$i]
}

# Move the selection to the newly-added final row, mirroring Excel's
# behaviour of leaving the cursor on the last entered cell.
$ws.Range("A44").Select()
